# Natmi following Dr Hou advice
# The LR-pair table (Bmp4 -> Rgmb) is recomputed with a new "ECs" sending/
# target cluster added alongside the existing FAPs / sCs clusters, so the
# former 2 (senders) x 3 (targets) = 6 data rows become 3 x 4 = 12 rows.
# Every numeric column (E:T) is recalculated for the new cluster set, and
# the "ECs" label is introduced as a new shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E..T the 16 numeric metric columns.
$data = @(
    ,("ECs",  "Bmp4", "Rgmb", "ECs",  2, 0.6666666666666666, 4.504599666666667,  13.513799,          0.2870666703033801, 0.2870666703033801, 2, 0.6666666666666666, 9.249291999999999,  27.747876,   0.1285586120341026, 0.1285586120341026, 41.66435766010266,  374.979218940924,  0.03690489269545389, 0.03690489269545388)
    ,("ECs",  "Bmp4", "Rgmb", "FAPs", 2, 0.6666666666666666, 4.504599666666667,  13.513799,          0.2870666703033801, 0.2870666703033801, 3, 1,                  35.42317866666667,  106.269536,  0.492357110492641,  0.492357110492641,  159.5672388141404,  1436.105149327264, 0.1413393163093158,  0.1413393163093159)
    ,("ECs",  "Bmp4", "Rgmb", "M2",   2, 0.6666666666666666, 4.504599666666667,  13.513799,          0.2870666703033801, 0.2870666703033801, 3, 1,                  3.703826,            11.111478,   0.0514805597850973, 0.05148055978509729,16.68425336499133,  150.158280284922,  0.01477835288286197, 0.01477835288286197)
    ,("ECs",  "Bmp4", "Rgmb", "sCs",  2, 0.6666666666666666, 4.504599666666667,  13.513799,          0.2870666703033801, 0.2870666703033801, 3, 1,                  23.569813,           70.709439,   0.3276037176881591, 0.3276037176881591, 106.1725717831957,  955.5531460487611, 0.09404410841574837, 0.09404410841574838)
    ,("FAPs", "Bmp4", "Rgmb", "ECs",  3, 1,                   7.480073,           22.440219,          0.476686011772755,  0.476686011772755,  2, 0.6666666666666666, 9.249291999999999,  27.747876,   0.1285586120341026, 0.1285586120341026, 69.18537935831598,  622.6684142248439, 0.06128209204957729, 0.06128209204957728)
    ,("FAPs", "Bmp4", "Rgmb", "FAPs", 3, 1,                   7.480073,           22.440219,          0.476686011772755,  0.476686011772755,  3, 1,                  35.42317866666667,  106.269536,  0.492357110492641,  0.492357110492641,  264.9679623187093,  2384.711660868384,0.2346997473686947,  0.2346997473686947)
    ,("FAPs", "Bmp4", "Rgmb", "M2",   3, 1,                   7.480073,           22.440219,          0.476686011772755,  0.476686011772755,  3, 1,                  3.703826,            11.111478,   0.0514805597850973, 0.05148055978509729,27.704888859298,    249.343999733682,  0.02454006272778691, 0.02454006272778691)
    ,("FAPs", "Bmp4", "Rgmb", "sCs",  3, 1,                   7.480073,           22.440219,          0.476686011772755,  0.476686011772755,  3, 1,                  23.569813,           70.709439,   0.3276037176881591, 0.3276037176881591, 176.303921836349,   1586.735296527141,0.1561641096266961,  0.1561641096266961)
    ,("sCs",  "Bmp4", "Rgmb", "ECs",  3, 1,                   3.707151333333333,  11.121454,          0.2362473179238649, 0.2362473179238649, 2, 0.6666666666666666, 9.249291999999999,  27.747876,   0.1285586120341026, 0.1285586120341026, 34.28852517018933,  308.596726531704,  0.03037162728907145, 0.03037162728907144)
    ,("sCs",  "Bmp4", "Rgmb", "FAPs", 3, 1,                   3.707151333333333,  11.121454,          0.2362473179238649, 0.2362473179238649, 3, 1,                  35.42317866666667,  106.269536,  0.492357110492641,  0.492357110492641,  131.3190840250382,  1181.871756225344,0.1163180468146304,  0.1163180468146304)
    ,("sCs",  "Bmp4", "Rgmb", "M2",   3, 1,                   3.707151333333333,  11.121454,          0.2362473179238649, 0.2362473179238649, 3, 1,                  3.703826,            11.111478,   0.0514805597850973, 0.05148055978509729,13.73064349433467,  123.575791449012,  0.01216214417444842, 0.01216214417444842)
    ,("sCs",  "Bmp4", "Rgmb", "sCs",  3, 1,                   3.707151333333333,  11.121454,          0.2362473179238649, 0.2362473179238649, 3, 1,                  23.569813,           70.709439,   0.3276037176881591, 0.3276037176881591, 87.37686368936734,  786.391773204306,  0.0773954996457146,  0.07739549964571461)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
